$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = newPriceText; E = newVolumeText } (only columns that changed)
$updates = @{
  2 = @{ D="34.883.57"; E="  -2.53%  " }
  3 = @{ D="1.809.04"; E="  -3.42%  " }
  4 = @{ E="  +0.11%  " }
  5 = @{ D="231.63"; E="  -0.52%  " }
  6 = @{ E="  -1.62%  " }
  7 = @{ E="  +0.13%  " }
  8 = @{ D="39.28"; E="  -7.60%  " }
  9 = @{ E="  +1.82%  " }
  10 = @{ D="0.0679"; E="  -3.03%  " }
  11 = @{ D="0.0995"; E="  -1.94%  " }
  12 = @{ D="2.070.54"; E="  -3.51%  " }
  13 = @{ D="1.840.14"; E="  -1.24%  " }
  14 = @{ E="  -4.21%  " }
  15 = @{ D="10.91"; E="  -7.03%  " }
  16 = @{ E="  -5.12%  " }
  17 = @{ D="34.825.02"; E="  -2.78%  " }
  18 = @{ D="69.05"; E="  -2.59%  " }
  19 = @{ D="0.0₃0780"; E="  -3.76%  " }
  20 = @{ E="  -4.81%  " }
  21 = @{ D="11.79"; E="  -5.89%  " }
  22 = @{ D="4.61"; E="  -4.60%  " }
  23 = @{ E="  +0.32%  " }
  24 = @{ E="  -1.07%  " }
  25 = @{ D="172.50"; E="  +0.38%  " }
  26 = @{ D="7.78"; E="  -3.98%  " }
  27 = @{ D="17.20"; E="  -4.66%  " }
  28 = @{ E="  -3.59%  " }
  29 = @{ E="  +6.78%  " }
  30 = @{ E="  +0.14%  " }
  31 = @{ E="  -0.49%  " }
  32 = @{ E="  -1.81%  " }
  33 = @{ E="  -4.01%  " }
  34 = @{ D="1.75"; E="  -8.84%  " }
  35 = @{ E="  +4.96%  " }
  36 = @{ D="0.682"; E="  -1.70%  " }
  37 = @{ D="91.00"; E="  -7.71%  " }
  38 = @{ D="1.30"; E="  +3.18%  " }
  39 = @{ D="1.306.32"; E="  -4.64%  " }
  40 = @{ E="  -3.72%  " }
  41 = @{ D="2.46"; E="  -1.42%  " }
  42 = @{ D="0.957"; E="  -7.13%  " }
  43 = @{ D="14.28"; E="  -5.59%  " }
  44 = @{ E="  -13.40%  " }
  45 = @{ E="  -5.35%  " }
  46 = @{ E="  -3.59%  " }
  48 = @{ D="1.991.79"; E="  -2.56%  " }
  49 = @{ D="0.0672"; E="  +6.90%  " }
  50 = @{ E="  +0.15%  " }
  51 = @{ D="98.72"; E="  -6.47%  " }
}

foreach ($row in $updates.Keys) {
  $rowData = $updates[$row]
  if ($rowData.ContainsKey("D")) {
    $cell = $ws.Cells.Item($row, 4)   # column D = Price
    $cell.NumberFormat = "@"          # force text storage so "231.63"-style values
                                       # don't get reinterpreted as numbers (which
                                       # would also strip meaningful trailing zeros,
                                       # e.g. "172.50" -> 172.5)
    $cell.Value = $rowData["D"]
    $cell.Style = "Normal"            # restore original (unstyled, General-format) cell style
  }
  if ($rowData.ContainsKey("E")) {
    $ws.Cells.Item($row, 5).Value = $rowData["E"]   # column E = Volume(1h); stays text naturally
  }
}
